# PA_Cadastrar_Aluno.docx edit
#
# 1) "Pre-condicoes" section: the sentence "Aluno  nao cadastrado" had a
#    double space wrapped in proofErr gramStart/gramEnd markers; collapse it
#    to a single space ("Aluno nao cadastrado") and drop the stray proofErr
#    bookkeeping runs.
# 2) Word's internal "_GoBack" bookmark (tracks the last edit position) is
#    relocated from the "[FA1]" alternate-flow paragraph ("O aluno faz as
#    alteracoes...") to the spot right after "Aluno" in the sentence that
#    was just edited above - exactly what real Word does automatically when
#    you type at a new location and save.

$d = $word.ActiveDocument

# --- 1) Fix "Aluno  não cadastrado" -> "Aluno não cadastrado" -------------
$rng = $d.Content
$rng.Find.Execute("Aluno  não cadastrado", $false, $false, $false, $false, $false, $true, 1, $false, "Aluno não cadastrado", 2)

# --- 2) Move the "_GoBack" bookmark to right after "Aluno" ----------------
$rng2 = $d.Content
$rng2.Find.Execute("Aluno não cadastrado", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$bmStart = $rng2.Start + 5
$bmRange = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRange)
